# Weekly update: insert the newest week's record for
# "Terminal La Palmera de La Serena - Perejil" as a new row 83,
# shifting the existing historical rows (83:106) down to (84:107).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83 (shifts rows 83-106 down to 84-107).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with this week's data. Everything
# besides the date (D) and the volume (J) repeats the values that were
# already in the old row 83 (now row 84).
$ws.Cells.Item(83, 1).Value = 8
$ws.Cells.Item(83, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = 44551
$ws.Cells.Item(83, 4).NumberFormat = $ws.Cells.Item(84, 4).NumberFormat
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = 100112044
$ws.Cells.Item(83, 7).Value = "Perejil"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 2880
$ws.Cells.Item(83, 11).Value = 2000
$ws.Cells.Item(83, 12).Value = 2500
$ws.Cells.Item(83, 13).Value = 2250
$ws.Cells.Item(83, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(83, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(83, 16).Value = 1500
$ws.Cells.Item(83, 17).Value = 1.5
$ws.Cells.Item(83, 18).Value = "Hortaliza"
